# Update Sheets via scheduled runner: refresh market-price derived columns
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 9273.75
$ws.Range("I74").Value = 9273.75
$ws.Range("K74").Value = 9273.75
$ws.Range("M74").Value = -8337.75
$ws.Range("H77").Value = 9273.75
$ws.Range("I77").Value = 9273.75
$ws.Range("K77").Value = 46368.75
$ws.Range("M77").Value = -41688.75
$ws.Range("H98").Value = 3911.875
$ws.Range("I98").Value = 1715.8334
$ws.Range("K98").Value = 1715.8334
$ws.Range("M98").Value = -217.8334
$ws.Range("H112").Value = 2721
$ws.Range("J112").Value = 2853.3845
$ws.Range("L112").Value = 8560.1535
$ws.Range("N112").Value = -10776.1535
$ws.Range("H122").Value = 3911.875
$ws.Range("I122").Value = 1715.8334
$ws.Range("K122").Value = 5147.5002
$ws.Range("M122").Value = -2697.5002
$ws.Range("H132").Value = 1096.4706
$ws.Range("I132").Value = 1096.4706
$ws.Range("K132").Value = 3289.4118
$ws.Range("M132").Value = -759.4118000000003
$ws.Range("H135").Value = 300
$ws.Range("J135").Value = 300
$ws.Range("L135").Value = 2700
$ws.Range("N135").Value = -7770
$ws.Range("H137").Value = 1124.25
$ws.Range("I137").Value = 1124.25
$ws.Range("K137").Value = 3372.75
$ws.Range("M137").Value = -822.75
$ws.Range("H138").Value = 2401.3447
$ws.Range("I138").Value = 1449.75
$ws.Range("J138").Value = 2553.6
$ws.Range("K138").Value = 4349.25
$ws.Range("L138").Value = 7660.799999999999
$ws.Range("M138").Value = 790.75
$ws.Range("N138").Value = -17940.8
$ws.Range("H141").Value = 3095
$ws.Range("I141").Value = 3095
$ws.Range("K141").Value = 9285
$ws.Range("M141").Value = -4105

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5971.722
$ws.Range("I32").Value = 6317.1177
$ws.Range("K32").Value = 6317.1177
$ws.Range("M32").Value = -6030.1177
$ws.Range("H45").Value = 3032.3635
$ws.Range("I45").Value = 3187
$ws.Range("J45").Value = 1486
$ws.Range("K45").Value = 3187
$ws.Range("L45").Value = 1486
$ws.Range("M45").Value = -2810
$ws.Range("N45").Value = -2240
$ws.Range("H61").Value = 2015.5
$ws.Range("I61").Value = 2015.5
$ws.Range("K61").Value = 2015.5
$ws.Range("M61").Value = -1803.5
$ws.Range("H74").Value = 9726.833000000001
$ws.Range("I74").Value = 9702.091
$ws.Range("K74").Value = 9702.091
$ws.Range("M74").Value = -8828.091
$ws.Range("H77").Value = 9726.833000000001
$ws.Range("I77").Value = 9702.091
$ws.Range("K77").Value = 48510.455
$ws.Range("M77").Value = -44142.455
$ws.Range("H122").Value = 5441.875
$ws.Range("I122").Value = 4933.4287
$ws.Range("K122").Value = 14800.2861
$ws.Range("M122").Value = -12350.2861
$ws.Range("H132").Value = 3572.6428
$ws.Range("I132").Value = 3588.7144
$ws.Range("K132").Value = 10766.1432
$ws.Range("M132").Value = -8236.143199999999
$ws.Range("H136").Value = 2015.5
$ws.Range("I136").Value = 2015.5
$ws.Range("K136").Value = 6046.5
$ws.Range("M136").Value = -3496.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3032.8333
$ws.Range("I134").Value = 3032.8333
$ws.Range("K134").Value = 9098.499899999999
$ws.Range("M134").Value = -6563.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2000.375
$ws.Range("I31").Value = 1691.4
$ws.Range("J31").Value = 2515.3333
$ws.Range("K31").Value = 1691.4
$ws.Range("L31").Value = 2515.3333
$ws.Range("M31").Value = -1396.4
$ws.Range("N31").Value = -3105.3333
$ws.Range("H34").Value = 2000.375
$ws.Range("I34").Value = 1691.4
$ws.Range("J34").Value = 2515.3333
$ws.Range("K34").Value = 1691.4
$ws.Range("L34").Value = 2515.3333
$ws.Range("M34").Value = -1489.4
$ws.Range("N34").Value = -2919.3333
$ws.Range("H134").Value = 10585.9
$ws.Range("I134").Value = 11133.143
$ws.Range("J134").Value = 9309
$ws.Range("K134").Value = 33399.429
$ws.Range("L134").Value = 27927
$ws.Range("M134").Value = -30864.429
$ws.Range("N134").Value = -32997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1502931
$ws.Range("I8").Value = 1502931
$ws.Range("K8").Value = 4508793
$ws.Range("M8").Value = -4508654
$ws.Range("H87").Value = 9154.333000000001
$ws.Range("I87").Value = 7785.2
$ws.Range("K87").Value = 23355.6
$ws.Range("M87").Value = -22107.6
$ws.Range("H90").Value = 9154.333000000001
$ws.Range("I90").Value = 7785.2
$ws.Range("K90").Value = 70066.8
$ws.Range("M90").Value = -63826.8
$ws.Range("H129").Value = 2851.4285
$ws.Range("J129").Value = 3093.3333
$ws.Range("L129").Value = 9279.999899999999
$ws.Range("N129").Value = -19279.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5699.8
$ws.Range("I102").Value = 4902.75
$ws.Range("J102").Value = 8888
$ws.Range("K102").Value = 4902.75
$ws.Range("L102").Value = 8888
$ws.Range("M102").Value = -3280.75
$ws.Range("N102").Value = -12132
$ws.Range("H122").Value = 4895.6665
$ws.Range("I122").Value = 4895.6665
$ws.Range("K122").Value = 14686.9995
$ws.Range("M122").Value = -12236.9995
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4641.5386
$ws.Range("I46").Value = 1445
$ws.Range("J46").Value = 43000
$ws.Range("K46").Value = 1445
$ws.Range("L46").Value = 43000
$ws.Range("M46").Value = -1257
$ws.Range("N46").Value = -43376
$ws.Range("H68").Value = 2649.2727
$ws.Range("I68").Value = 2714.2
$ws.Range("K68").Value = 2714.2
$ws.Range("M68").Value = -1965.2
$ws.Range("H71").Value = 2649.2727
$ws.Range("I71").Value = 2714.2
$ws.Range("K71").Value = 13571
$ws.Range("M71").Value = -9827
$ws.Range("H132").Value = 5641.7144
$ws.Range("I132").Value = 4750
$ws.Range("K132").Value = 14250
$ws.Range("M132").Value = -11720
$ws.Range("H136").Value = 3717.4
$ws.Range("I136").Value = 3717.4
$ws.Range("K136").Value = 11152.2
$ws.Range("M136").Value = -8602.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2603.8
$ws.Range("I132").Value = 1507
$ws.Range("K132").Value = 4521
$ws.Range("M132").Value = -1991
$ws.Range("H136").Value = 8208.125
$ws.Range("I136").Value = 7333.2
$ws.Range("K136").Value = 21999.6
$ws.Range("M136").Value = -19449.6

Write-Host "Applied 165 cell updates across 8 sheets"